# "ajax kind of export" -- the sheet now reflects an import-code / material
# lookup instead of supplier/good ids.
#
# Shared-string pool order matters (diff expects IP00034 @3, code_import @4,
# id_material_detail @5) so write the new data value before the new headers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A's data rows become the text code "IP00034" instead of the
# numeric supplier id 1.
$ws.Range("A2").Value = "IP00034"
$ws.Range("A3").Value = "IP00034"

# Rename the header row: id_supplier -> code_import, id_good -> id_material_detail.
$ws.Range("A1").Value = "code_import"
$ws.Range("B1").Value = "id_material_detail"

# Column B's header text got a lot longer, so both columns get re-sized to fit.
$ws.Columns("A").ColumnWidth = 10.666666666666666
$ws.Columns("B").ColumnWidth = 15

# Leave the selection where the author last clicked.
[void]$ws.Range("E9").Select()
